# Suivi.xlsx - "feat: code vhdl in ROM"
# Fill in the two new work-log entries (row 15 and row 16) on the "Travail"
# sheet, mirroring the left block (B:D) into the right block (F:H), and
# update the current selection to F13.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Travail")

# ---- Row 15 : 14.05.2024 - 1h30 - "Création de la structure hiérarchique" ----

# Left block
$ws.Range("B14").Copy()
$ws.Range("B15").PasteSpecial(-4122)
$ws.Range("B15").Value = 45426

$ws.Range("C14").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("C15").Value = "1h30"

$ws.Range("D14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = "Création de la structure hiérarchique"

# Right block (mirrors the left block)
$ws.Range("F14").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("F15").Value = 45426

$ws.Range("G14").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("G15").Value = "1h30"

$ws.Range("H14").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("H15").Value = "Création de la structure hiérarchique"

# ---- Row 16 : 19.05.2024 - 3h00 - "Finission de la structure hiérarchique et création du code dans la ROM" ----

# Left block (Travail entered before Durée, matching shared-string order)
$ws.Range("B13").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("B16").Value = 45431

$ws.Range("D13").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("D16").Value = "Finission de la structure hiérarchique et création du code dans la ROM"

$ws.Range("C13").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C16").Value = "3h00"

# Right block (mirrors the left block)
$ws.Range("F13").Copy()
$ws.Range("F16").PasteSpecial(-4122)
$ws.Range("F16").Value = 45431

$ws.Range("H13").Copy()
$ws.Range("H16").PasteSpecial(-4122)
$ws.Range("H16").Value = "Finission de la structure hiérarchique et création du code dans la ROM"

$ws.Range("G13").Copy()
$ws.Range("G16").PasteSpecial(-4122)
$ws.Range("G16").Value = "3h00"

# Row 16 grows to a 2-line row because of the long wrapped text above.
$ws.Rows.Item(16).RowHeight = 28.8

$excel.CutCopyMode = 0

# Move the active selection like the author's last click before saving.
$ws.Range("F13").Select()
